# Generate Report for Handback
#
# Mirrors the "Latest Target File" (F) and "Latest Handback File" (G)
# columns onto rows 2 & 3 of the zh-cn / de-de sheets (the handback
# arrived in sync with the original handoff, so these reuse the same
# .md / .xlf targets already hyperlinked from columns A / D), stamps
# the handback datetime into column H, and flips the Overview sheet's
# status text from "Ready for handoff" to "Handed back: in sync with
# en-US".

$wb = $excel.ActiveWorkbook

# Blue/underline color used by the workbook's existing hyperlink cells
# (A2/A3/D2/D3), expressed as an OLE (BGR) color so new cells visually
# match instead of falling back to the theme default.
$hyperlinkColor = 15570276   # RGB(0x64,0x95,0xED) -> R + G*256 + B*65536

# ---------------------------------------------------------------------
# Overview sheet: "Ready for handoff" -> "Handed back: in sync with en-US"
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$newStatus = "Handed back: in sync with en-US"
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# Helper: add target-file (F) / handback-file (G) hyperlink cells for a
# handoff/handback pair of rows on a language sheet, plus the handback
# datetime stamp in column H.
#
# NB: positional parameters only -- this host's PowerShell-subset
# function binder does not reliably bind `-named` arguments.
# ---------------------------------------------------------------------
function Set-HandbackRow($ws, $row, $mdDisplay, $mdUrl, $xlfDisplay, $xlfUrl, $handbackDateTime) {
    $fCell = $ws.Range("F" + $row)
    $fCell.Value = $mdDisplay
    $ws.Hyperlinks.Add($fCell, $mdUrl, "", "", $mdDisplay)
    $fCell.Font.Underline = 2
    $fCell.Font.Color = $hyperlinkColor

    $gCell = $ws.Range("G" + $row)
    $gCell.Value = $xlfDisplay
    $ws.Hyperlinks.Add($gCell, $xlfUrl, "", "", $xlfDisplay)
    $gCell.Font.Underline = 2
    $gCell.Font.Color = $hyperlinkColor

    $ws.Range("H" + $row).Value = $handbackDateTime
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

Set-HandbackRow $zhcn 2 `
    "0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/29711e7fdf5ec62afc8d72259ecba2f945528085/e2e/0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.md" `
    "0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.7d6e6f73e86ad0afe8eb81a66213da2cba5f59c6.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/62c9f873c35295e4459970c75ca8ba3939a48afe/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.7d6e6f73e86ad0afe8eb81a66213da2cba5f59c6.zh-cn.xlf" `
    "2016-03-21 10:43:09"

Set-HandbackRow $zhcn 3 `
    "ff945447-2701-4791-9b2d-41e05a4160a7.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/29711e7fdf5ec62afc8d72259ecba2f945528085/e2e/ff945447-2701-4791-9b2d-41e05a4160a7.md" `
    "ff945447-2701-4791-9b2d-41e05a4160a7.93c8d1b8c9ab71c0f1001f53f2d8d805e7efbae6.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/62c9f873c35295e4459970c75ca8ba3939a48afe/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ff945447-2701-4791-9b2d-41e05a4160a7.93c8d1b8c9ab71c0f1001f53f2d8d805e7efbae6.zh-cn.xlf" `
    "2016-03-21 10:43:09"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

Set-HandbackRow $dede 2 `
    "0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/29711e7fdf5ec62afc8d72259ecba2f945528085/e2e/0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.md" `
    "0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.7d6e6f73e86ad0afe8eb81a66213da2cba5f59c6.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8410f360bd6907985443a2555ca5e17d37425994/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.7d6e6f73e86ad0afe8eb81a66213da2cba5f59c6.de-de.xlf" `
    "2016-03-21 10:43:17"

Set-HandbackRow $dede 3 `
    "ff945447-2701-4791-9b2d-41e05a4160a7.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/29711e7fdf5ec62afc8d72259ecba2f945528085/e2e/ff945447-2701-4791-9b2d-41e05a4160a7.md" `
    "ff945447-2701-4791-9b2d-41e05a4160a7.93c8d1b8c9ab71c0f1001f53f2d8d805e7efbae6.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8410f360bd6907985443a2555ca5e17d37425994/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ff945447-2701-4791-9b2d-41e05a4160a7.93c8d1b8c9ab71c0f1001f53f2d8d805e7efbae6.de-de.xlf" `
    "2016-03-21 10:43:17"

Write-Output "Handback report generated."
